# Build site at 2021-10-27 12:19:33 UTC
#
# LOQ4247.xlsx changes:
#  1. "Semestre ideal:" value changes from "EP-10" to "EA-1,EP-10"
#     (both the B and C columns of that row share the same text).
#  2. The trailing "Requisitos:" row and its associated requirement text
#     row are removed entirely (rows 23 and 24), shrinking the used range
#     from A1:C24 down to A1:C22.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update "Semestre ideal:" value (row 9, columns B and C) in place.
$ws.Range("B9:C9").Value = "EA-1,EP-10"

# 2. Remove the last two rows (Requisitos: / LOB1009 ...) completely,
#    shifting everything below them up (there is nothing below, so this
#    just shrinks the sheet's dimension).
$ws.Range("A23:C24").EntireRow.Delete()
